$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: value 2.5 moves from column C to column B
$ws.Range("C2").ClearContents()
$ws.Range("B2").Value = 2.5

# Row 5: B5 becomes 5.5, C5 cleared (was 2)
$ws.Range("C5").ClearContents()
$ws.Range("B5").Value = 5.5

# Row 31: B31 2.5 -> 2.75, C31 1.25 -> 1
$ws.Range("B31").Value = 2.75
$ws.Range("C31").Value = 1

# Row 38: B38 becomes 3, C38 cleared (was 0.75)
$ws.Range("C38").ClearContents()
$ws.Range("B38").Value = 3

# Row 54: B54 becomes 3.25, C54 cleared (was 2.5)
$ws.Range("C54").ClearContents()
$ws.Range("B54").Value = 3.25

# Update sheet view: remove frozen/scrolled topLeftCell, reset selection to B5
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
